$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 29-33: time value in column A (same time number-format as the
# existing A3:A28 entries), course topic text in column B (new shared
# strings).
$times = @(
    0.62847222222222221,
    0.63194444444444442,
    0.64930555555555558,
    0.66875000000000007,
    0.68055555555555547
)
$notes = @(
    "SQL Injection 寫惡意語法 偷窺 破壞資料 原因: 寫SQL時串接語法",
    "SQL DateTime",
    "百元買百雞 練習作業",
    "找練習題",
    "防止SQL Injection攻擊"
)

$timeFormat = $ws.Range("A13").NumberFormat

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = 29 + $i
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $times[$i]
    $aCell.NumberFormat = $timeFormat

    $ws.Cells.Item($row, 2).Value = $notes[$i]
}

# Move the view: scroll so row 22 is near the top and select B34, matching
# the end-of-session cursor position recorded in the workbook.
$ws.Range("A22").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("B34").Select() | Out-Null
